$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'23.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.237"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05813"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'6.465"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'3.224"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8081"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8855"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1398"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07135"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03102"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03046"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09333"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.829"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001555"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04728"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0006024"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'17OneONE"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.006202"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.001261"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.004070"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.00008705"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'3.540"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").Value = "'0.1322"
$ws.Range("D26").Style = "Normal"
$ws.Range("D40").Value = "'0.03854"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006294"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'40KickTokenKICKBestin24h"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1053"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.002515"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.007841"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005330"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D46").Style = "Normal"
$ws.Range("E47").Value = "'46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002835"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'47BOLOBOLO"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D49").Style = "Normal"
